# The workbook's data table started at row 2 (row 1 was a blank leading
# row). This edit removes that blank row, shifting the whole table -- and
# the floating "TextBox 1" callout that is anchored relative to the grid --
# up by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1").Delete()

# The text box's cell-relative anchor needs to follow the row shift (Excel
# keeps shapes pinned to the grid when rows move). Nudge it up by one
# row's worth of points so it tracks the same content it annotated before.
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top - $ws.Rows("1").Height

# Leave the selection where the author's session ended up after the edit.
$ws.Range("B10").Select() | Out-Null
